$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos data refresh: update price (D) and volume-1h (E) columns for all
# rows, plus coin name (B) and link (C) for rows whose ranking position
# shifted (WrappedeETH dropped out of the top list, dogwifhat entered it).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.559.80'
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.348.05'
$ws.Range("E3").Value = '  -2.46%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.56'
$ws.Range("E5").Value = '  -3.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '656.24'
$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("E7").Value = '  -6.76%  '

$ws.Range("E8").Value = '  -5.71%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.995'
$ws.Range("E10").Value = '  -9.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.343.30'
$ws.Range("E11").Value = '  -2.55%  '

$ws.Range("E12").Value = '  -3.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.81'
$ws.Range("E13").Value = '  -4.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '97.262.37'
$ws.Range("E14").Value = '  -0.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.07'
$ws.Range("E15").Value = '  -5.21%  '

$ws.Range("E16").Value = '  -6.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.980.50'
$ws.Range("E17").Value = '  -2.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.62'
$ws.Range("E18").Value = '  -7.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.361.45'
$ws.Range("E19").Value = '  -2.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.50'
$ws.Range("E20").Value = '  -2.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.511'
$ws.Range("E21").Value = '  -16.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.76'
$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '506.32'
$ws.Range("E23").Value = '  -2.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.33'
$ws.Range("E24").Value = '  -5.71%  '

$ws.Range("E25").Value = '  -4.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.81'
$ws.Range("E26").Value = '  +6.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.82'
$ws.Range("E27").Value = '  -6.68%  '

$ws.Range("E28").Value = '  -8.49%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.29'
$ws.Range("E29").Value = '  -5.11%  '

$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.140'
$ws.Range("E31").Value = '  -10.70%  '

$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.185'
$ws.Range("E32").Value = '  -8.05%  '

$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.55'
$ws.Range("E33").Value = '  +8.07%  '

$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.35%  '

$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.554'
$ws.Range("E35").Value = '  -6.60%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.33'
$ws.Range("E36").Value = '  -6.56%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.94'
$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.51'
$ws.Range("E38").Value = '  +2.87%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '524.59'
$ws.Range("E39").Value = '  -2.00%  '

$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.151'
$ws.Range("E41").Value = '  -3.23%  '

$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.39'
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("B43").Value = 'ImmutableX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.73'
$ws.Range("E43").Value = '  +2.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.839'
$ws.Range("E44").Value = '  -5.32%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0423'
$ws.Range("E45").Value = '  -5.25%  '

$ws.Range("E46").Value = '  +7.30%  '

$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.53'
$ws.Range("E47").Value = '  -7.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.58'
$ws.Range("E48").Value = '  -4.01%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.41'
$ws.Range("E49").Value = '  -9.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.60'
$ws.Range("E50").Value = '  +6.05%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.15'
$ws.Range("E51").Value = '  -6.12%  '
